# Applies the row-content rearrangement described in the commit diff.
# Rows 4<->5 swap, rows 21<->22 swap, and rows 24/25/26 rotate (24<-26, 25<-24, 26<-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 130873745
$ws.Range("B4").Value = 79243
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = 'Garnlav'
$ws.Range("G4").Value = 'Alectoria sarmentosa'
$ws.Range("H4").Value = '(Ach.) Ach.'
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 438633
$ws.Range("R4").Value = 6795187
$ws.Range("AX4").Value = 'Eva Löfqvist'
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("N4").ClearContents()

# Row 5
$ws.Range("A5").Value = 130873700
$ws.Range("B5").Value = 57076
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 102613
$ws.Range("F5").Value = 'Orre'
$ws.Range("G5").Value = 'Lyrurus tetrix'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("M5").Value = 'färska spår'
$ws.Range("Q5").Value = 438768
$ws.Range("R5").Value = 6795206
$ws.Range("AX5").Value = 'Eva Löfqvist, Alfhild Sehlin'

# Row 21
$ws.Range("A21").Value = 130873726
$ws.Range("B21").Value = 79243
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("M21").ClearContents()
$ws.Range("Q21").Value = 438662
$ws.Range("R21").Value = 6795157
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("N21").ClearContents()

# Row 22
$ws.Range("A22").Value = 130873699
$ws.Range("B22").Value = 57076
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 102613
$ws.Range("F22").Value = 'Orre'
$ws.Range("G22").Value = 'Lyrurus tetrix'
$ws.Range("H22").Value = '(Linnaeus, 1758)'
$ws.Range("M22").Value = 'färska spår'
$ws.Range("Q22").Value = 438808
$ws.Range("R22").Value = 6795184

# Row 24
$ws.Range("A24").Value = 130873702
$ws.Range("B24").Value = 8451
$ws.Range("D24").Value = 'LC'
$ws.Range("E24").Value = 106545
$ws.Range("F24").Value = 'Mindre märgborre'
$ws.Range("G24").Value = 'Tomicus minor'
$ws.Range("H24").Value = '(Hartig, 1834)'
$ws.Range("M24").Value = 'äldre gnagspår'
$ws.Range("Q24").Value = 438985
$ws.Range("R24").Value = 6795128

# Row 25
$ws.Range("A25").Value = 130873724
$ws.Range("Q25").Value = 438724
$ws.Range("R25").Value = 6795185

# Row 26
$ws.Range("A26").Value = 130873735
$ws.Range("B26").Value = 79243
$ws.Range("D26").Value = 'NT'
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = 'Garnlav'
$ws.Range("G26").Value = 'Alectoria sarmentosa'
$ws.Range("H26").Value = '(Ach.) Ach.'
$ws.Range("M26").ClearContents()
$ws.Range("Q26").Value = 438776
$ws.Range("R26").Value = 6795246
$ws.Range("J26").ClearContents()
$ws.Range("K26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("AF26").ClearContents()
